$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "35.287.87"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.24%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "1.910.93"
$ws.Range("D3").Style = "Normal"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.722"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +9.23%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'254.57"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.55%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'40.64"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.72%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "  +5.27%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'52.73"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.14%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.0763"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +6.19%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.0986"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.56%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "2.187.27"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.02%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'12.76"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.23%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("E15").Value = "  +3.14%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'4.94"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.35%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "1.919.61"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.38%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "35.281.71"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.18%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "  +2.37%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "0.0₃0853"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.59%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'243.90"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.68%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("E23").Value = "  +5.73%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "  +0.13%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "  +4.90%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "  +3.78%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'167.18"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.61%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'8.67"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.37%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'18.78"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.69%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = "  +4.76%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "4.128.89"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +19.46%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'4.37"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +5.27%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = "  +14.24%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = "  +23.66%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = "  +3.73%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = "  +2.63%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = "  +0.05%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'0.913"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.85%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'2.05"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.28%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = "  +5.01%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'17.09"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.12%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'96.81"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +7.51%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "  +0.92%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "  +2.24%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "1.337.30"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.19%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "  +1.29%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "  +0.93%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("E49").Value = "  -0.73%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'45.28"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -5.59%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'11.83"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +12.94%  "
$ws.Range("E51").Style = "Normal"
